$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.140.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.88%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.862.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.54%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.17%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'234.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.28%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.17%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4667"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.2826"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.77%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.30%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'20.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.16%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07820"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.21%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'95.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -7.59%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.857.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.01%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.122"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.88%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.6700"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.51%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'279.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -5.08%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'30.183.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.94%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +0.09%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'5.447"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.83%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.51%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'2.098.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.52%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.000007239"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -4.50%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D24").Value = "'6.145"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.85%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'9.313"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.20%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'165.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.41%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'18.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -4.21%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.901"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -8.88%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.345"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.20%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.09576"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.76%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.397"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.42%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.468"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.27%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.118"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -4.88%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.04659"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.55%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.7005"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.85%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.096"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.92%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.705"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.88%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01859"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -5.20%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'6.278"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -8.71%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.527"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.82%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'72.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.81%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.8524"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.30%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.926"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -5.23%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.18%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.4156"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.60%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'103.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.56%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'988.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.49%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'7.118"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -5.58%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'9.132"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.53%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'33.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.43%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.1137"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.64%  "
$ws.Range("E51").Style = "Normal"
